$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for "Documents / md" right before the existing
#    "Documents / html" row (currently row 18).
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "Documents"
$ws.Range("B18").Value = "md"

# 2) Insert a new row for "Videos / gif" right after "Videos / mp4"
#    (which, after the previous insert, now sits on row 23).
$ws.Rows.Item(24).Insert()
$ws.Range("A24").Value = "Videos"
$ws.Range("B24").Value = "gif"

# 3) Remove the "Fonts / ttf" row (now row 26).
$ws.Rows.Item(26).Delete()

# 4) Swap the order of "Programming / json" and "Programming / py"
#    (now rows 29 and 30) so "py" comes first.
$ws.Range("B29").Value = "py"
$ws.Range("B30").Value = "json"

# The table ("ConfigTable") needs to grow by one row overall (32 -> 33
# data rows, i.e. range A1:B33 -> A1:B34) to track the net row change.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B34"))

# Reflect the scrolled viewport from the saved workbook (the sheet was
# left scrolled down near row 13 rather than showing the old L7 selection).
$ws.Range("A13").Select() | Out-Null

Write-Host "Dim after edits:" $ws.UsedRange.Address()
Write-Host "Table range:" $tbl.Range.Address()
